$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values (rows 2-6, columns A:C) per the "bs env for sicmdp" change
$arr = New-Object 'object[,]' 5,3
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 1.5
$arr[1,0] = 200
$arr[1,1] = 0
$arr[1,2] = 1.5
$arr[2,0] = 0
$arr[2,1] = 200
$arr[2,2] = 1.5
$arr[3,0] = 200
$arr[3,1] = 200
$arr[3,2] = 1.5
$arr[4,0] = 100
$arr[4,1] = 100
$arr[4,2] = 1.5
$ws.Range("A2:C6").Value = $arr

# Update the selected cell in the sheet view
$ws.Range("E9").Select()

# Update the workbook window height
$excel.ActiveWindow.Height = 17655
